$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.103.67"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.901.29"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'566.05"
$ws.Range("E5").Value = "  -3.51%  "
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.500"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "2.897.75"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").Value = "'6.93"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").Value = "'0.148"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").Value = "'0.0000237"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "'32.53"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "3.379.78"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "62.012.04"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").Value = "2.897.79"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").Value = "'427.87"
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").Value = "'13.04"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("D24").Value = "'78.40"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("D25").Value = "'11.96"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D27").Value = "'9.82"
$ws.Range("E27").Value = "  -3.77%  "
$ws.Range("E28").Value = "  -3.82%  "
$ws.Range("E29").Value = "  +4.36%  "
$ws.Range("D30").Value = "'6.91"
$ws.Range("E30").Value = "  -4.14%  "
$ws.Range("D31").Value = "'2.48"
$ws.Range("E31").Value = "  -3.25%  "
$ws.Range("E32").Value = "  -5.33%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "'25.57"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("D36").Value = "'0.954"
$ws.Range("E36").Value = "  -2.02%  "
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("D38").Value = "'48.80"
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").Value = "'2.92"
$ws.Range("E39").Value = "  -4.77%  "
$ws.Range("E40").Value = "  -6.38%  "
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("D42").Value = "'40.69"
$ws.Range("E42").Value = "  +4.75%  "
$ws.Range("D43").Value = "'8.05"
$ws.Range("E43").Value = "  -3.37%  "
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("D45").Value = "2.700.63"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "'133.58"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").Value = "'345.22"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("D50").Value = "'0.000217"
$ws.Range("E50").Value = "  +13.93%  "
$ws.Range("E51").Value = "  -1.56%  "
